$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data point was inserted at row 214, pushing all the existing
# rows (214..283) down by one (215..284). Insert a new row first so the
# rest of the sheet (and its formatting) shifts down automatically.
$ws.Rows.Item(214).Insert()

# Populate the newly inserted row 214 with the new reading. Columns that are
# constant across every row in this sheet (market/region/category/etc.) are
# copied from the surrounding rows; D (Fecha), J (Volumen) and O (Origen)
# carry the new values.
$ws.Range("A214").Value = 10
$ws.Range("B214").Value = "Vega Modelo de Temuco"
$ws.Range("C214").Value = "La Araucanía"
$ws.Range("D214").Value = 44524
$ws.Range("E214").Value = 9
$ws.Range("F214").Value = 100112008
$ws.Range("G214").Value = "Coliflor"
$ws.Range("H214").Value = "Sin especificar"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 1400
$ws.Range("K214").Value = 900
$ws.Range("L214").Value = 900
$ws.Range("M214").Value = 900
$ws.Range("N214").Value = "`$/unidad"
$ws.Range("O214").Value = "Región de O'Higgins"
$ws.Range("P214").Value = 900
$ws.Range("Q214").Value = 1
$ws.Range("R214").Value = "Hortaliza"
